$d = $word.ActiveDocument

# 1) "Convênio ou Contrato nº 001/2021" -> "Convênio nº 001/2021"
#    Remove the words "ou Contrato " so the remaining space run (spacing -3)
#    that used to separate "Convênio" and "ou" is the one left in place.
$d.Content.Find.Execute("ou Contrato ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "", 2)

# 2) "Data:27/{mes}/2021" -> "Data: {final}"
#    Replace the date placeholder with " {final}" built from four distinct
#    runs (" ", "{", "final", "}") so formatting/run layout matches the
#    target structure, leaving the existing "Data" / ":" runs untouched.
$r = $d.Content
$r.Find.Execute("27/{mes}/2021", $true, $false, $false, $false, $false, `
                 $true, 1, $false, "", 2)
$r.InsertAfter(" ")
$r.Collapse(0)
$r.InsertAfter("{")
$r.Collapse(0)
$r.InsertAfter("final")
$r.Collapse(0)
$r.InsertAfter("}")
